$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "LP1912"  (columns: A=Fecha/meta, B=Hora_Scrap, C=Hora_Llegada,
#                      D=Linea, E=Minutos, F=Parada, G=Fecha)
# Add 14 new scraped rows (237-250) and refresh the header metadata.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 30/12/2025 12:58:42"
$ws1.Range("A3").Value = "Total filas: 249"

$sheet1Rows = @(
    @("", "12:58:31", "13:04", "23_HERNANDEZ",         6, "LP1912", "30/12/2025"),
    @("", "12:58:31", "13:06", "16_P MOR-SANTA ANA",   8, "LP1912", "30/12/2025"),
    @("", "12:58:31", "13:08", "10_OLMOS",            10, "LP1912", "30/12/2025"),
    @("", "12:58:31", "13:19", "10_OLMOS",            21, "LP1912", "30/12/2025"),
    @("", "12:58:31", "13:26", "14_ABASTO",           28, "LP1912", "30/12/2025"),
    @("", "12:58:31", "13:32", "10_OLMOS",            34, "LP1912", "30/12/2025"),
    @("", "12:58:31", "13:34", "23_HERNANDEZ",        36, "LP1912", "30/12/2025"),
    @("", "12:58:31", "13:36", "15_ABASTO",           38, "LP1912", "30/12/2025"),
    @("", "12:58:31", "13:46", "16_SANTA ANA",        48, "LP1912", "30/12/2025"),
    @("", "12:58:31", "13:46", "17_ROMERO",           48, "LP1912", "30/12/2025"),
    @("", "12:58:31", "13:56", "16_P MOR-167 Y 521",  58, "LP1912", "30/12/2025"),
    @("", "12:58:31", "14:04", "17_ROMERO",           66, "LP1912", "30/12/2025"),
    @("", "12:58:31", "14:07", "23_HERNANDEZ",        69, "LP1912", "30/12/2025"),
    @("", "12:58:31", "14:21", "26_HERNANDEZ",        83, "LP1912", "30/12/2025")
)

$r = 237
foreach ($row in $sheet1Rows) {
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $ws1.Cells.Item($r, 4).Value = $row[3]
    $ws1.Cells.Item($r, 5).Value = $row[4]
    $ws1.Cells.Item($r, 6).Value = $row[5]
    $ws1.Cells.Item($r, 7).Value = $row[6]
    $r++
}

# ---------------------------------------------------------------------------
# Sheet 2: "LP1912-215" -- only the "last updated" header timestamp changes.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: 30/12/2025 12:58:42"

# ---------------------------------------------------------------------------
# Sheet 3: "6203-6173" (columns: A=meta, B=Fecha, C=Hora_Scrap, D=Hora_Llegada,
#                        E=Linea, F=Minutos, G=Parada)
# Add 2 new scraped rows (35-36) and refresh the header metadata.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 30/12/2025 12:58:42"
$ws3.Range("A3").Value = "Total filas: 35"

$sheet3Rows = @(
    @("", "30/12/2025", "12:58:42", "13:30", "215B_LP-P MOR-1 Y 57", 32, "L6173"),
    @("", "30/12/2025", "12:58:42", "14:09", "215A_LA PLATA",        71, "L6173")
)

$r = 35
foreach ($row in $sheet3Rows) {
    $ws3.Cells.Item($r, 2).Value = $row[1]
    $ws3.Cells.Item($r, 3).Value = $row[2]
    $ws3.Cells.Item($r, 4).Value = $row[3]
    $ws3.Cells.Item($r, 5).Value = $row[4]
    $ws3.Cells.Item($r, 6).Value = $row[5]
    $ws3.Cells.Item($r, 7).Value = $row[6]
    $r++
}
